# Update ranking at 2025-12-03 20:10
# Appends a new row (row 37) to the ranking sheet with the latest
# snapshot timestamp and placeholder "-" values for the ranking columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 37

$ws.Cells.Item($newRow, 1).Value = "2025/12/04 05:00"
$ws.Cells.Item($newRow, 2).Value = "-"
$ws.Cells.Item($newRow, 3).Value = "-"
$ws.Cells.Item($newRow, 4).Value = "-"
$ws.Cells.Item($newRow, 5).Value = "-"
$ws.Cells.Item($newRow, 6).Value = "-"
$ws.Cells.Item($newRow, 7).Value = "-"
